$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.853016668497649
$ws.Range("D2").Value = 9.088369646221008
$ws.Range("E2").Value = 10.96689336586168
$ws.Range("F2").Value = 49.13524593267123
$ws.Range("G2").Value = 3.695012136940502
$ws.Range("K2").Value = 12.34625181915244
$ws.Range("M2").Value = 15.36039177750424

$ws.Range("B3").Value = 7.785858890209346
$ws.Range("D3").Value = 8.939182259227128
$ws.Range("E3").Value = 10.74600278831383
$ws.Range("F3").Value = 47.83945206849071
$ws.Range("G3").Value = 3.699445129012229
$ws.Range("K3").Value = 12.21652037861133
$ws.Range("M3").Value = 15.22834356188051

$ws.Range("B4").Value = 7.746300445866833
$ws.Range("D4").Value = 8.845951120845916
$ws.Range("E4").Value = 10.61058355509153
$ws.Range("F4").Value = 47.03086974883865
$ws.Range("G4").Value = 3.70230162081969
$ws.Range("K4").Value = 12.14370142880355
$ws.Range("M4").Value = 15.15283563973332

$ws.Range("B5").Value = 7.730618163224602
$ws.Range("D5").Value = 8.807572175548442
$ws.Range("E5").Value = 10.55552102700945
$ws.Range("F5").Value = 46.69850302701713
$ws.Range("G5").Value = 3.703499667748156
$ws.Range("K5").Value = 12.11578703011279
$ws.Range("M5").Value = 15.12349705277274

$ws.Range("B6").Value = 7.728041057909301
$ws.Range("D6").Value = 8.801176726863259
$ws.Range("E6").Value = 10.54638734644501
$ws.Range("F6").Value = 46.64315286032549
$ws.Range("G6").Value = 3.703700660721014
$ws.Range("K6").Value = 12.11125924825798
$ws.Range("M6").Value = 15.11871269778988

$ws.Range("B7").Value = 7.74608715407609
$ws.Range("D7").Value = 8.845435062809056
$ws.Range("E7").Value = 10.60984037651952
$ws.Range("F7").Value = 47.02639838790589
$ws.Range("G7").Value = 3.702317640217812
$ws.Range("K7").Value = 12.14331778958657
$ws.Range("M7").Value = 15.15243413486471

$ws.Range("B8").Value = 7.829522803779549
$ws.Range("D8").Value = 9.037281798810435
$ws.Range("E8").Value = 10.89072839526425
$ws.Range("F8").Value = 48.69137679799665
$ws.Range("G8").Value = 3.696512785217299
$ws.Range("K8").Value = 12.30012751948261
$ws.Range("M8").Value = 15.31372494854361

$ws.Range("B9").Value = 8.005682116850684
$ws.Range("D9").Value = 9.399577020362591
$ws.Range("E9").Value = 11.44021514357754
$ws.Range("F9").Value = 51.83731866180889
$ws.Range("G9").Value = 3.686190615225303
$ws.Range("K9").Value = 12.65983502525913
$ws.Range("M9").Value = 15.67270317351568

$ws.Range("B10").Value = 8.141687567077364
$ws.Range("D10").Value = 9.656007478703019
$ws.Range("E10").Value = 11.83924348969237
$ws.Range("F10").Value = 54.05601213973218
$ws.Range("G10").Value = 3.679243999468526
$ws.Range("K10").Value = 12.9528985653224
$ws.Range("M10").Value = 15.96025020760918

$ws.Range("B11").Value = 8.204752672120771
$ws.Range("D11").Value = 9.770296379178843
$ws.Range("E11").Value = 12.01898669556042
$ws.Range("F11").Value = 55.04150425789497
$ws.Range("G11").Value = 3.676220030211075
$ws.Range("K11").Value = 13.0917575491844
$ws.Range("M11").Value = 16.09571267144426

$ws.Range("B12").Value = 8.228785462754319
$ws.Range("D12").Value = 9.813214920230823
$ws.Range("E12").Value = 12.08673519733691
$ws.Range("F12").Value = 55.4109839321648
$ws.Range("G12").Value = 3.675094335902972
$ws.Range("K12").Value = 13.14507426870945
$ws.Range("M12").Value = 16.14763231075602

$ws.Range("B13").Value = 8.22360316521222
$ws.Range("D13").Value = 9.803987978618796
$ws.Range("E13").Value = 12.07215934734005
$ws.Range("F13").Value = 55.33157886905187
$ws.Range("G13").Value = 3.675335913133483
$ws.Range("K13").Value = 13.1335599399529
$ws.Range("M13").Value = 16.1364235349755

$ws.Range("B14").Value = 8.206726947764773
$ws.Range("D14").Value = 9.773834645523333
$ws.Range("E14").Value = 12.02456708200107
$ws.Range("F14").Value = 55.07197734887568
$ws.Range("G14").Value = 3.676127030379056
$ws.Range("K14").Value = 13.09612960676029
$ws.Range("M14").Value = 16.09997191737029

$ws.Range("B15").Value = 8.196408883935575
$ws.Range("D15").Value = 9.755317362112399
$ws.Range("E15").Value = 11.99537249987253
$ws.Range("F15").Value = 54.91247350670439
$ws.Range("G15").Value = 3.676614136742439
$ws.Range("K15").Value = 13.07329611787685
$ws.Range("M15").Value = 16.07772397346962

$ws.Range("B16").Value = 8.137588420805992
$ws.Range("D16").Value = 9.648489117944225
$ws.Range("E16").Value = 11.82745581575553
$ws.Range("F16").Value = 53.99110425914664
$ws.Range("G16").Value = 3.67944434835928
$ws.Range("K16").Value = 12.94392996530854
$ws.Range("M16").Value = 15.95148715181112

$ws.Range("B17").Value = 8.1017956277837
$ws.Range("D17").Value = 9.5823334664802
$ws.Range("E17").Value = 11.72394492480835
$ws.Range("F17").Value = 53.41956896373213
$ws.Range("G17").Value = 3.681215337695507
$ws.Range("K17").Value = 12.86594312991125
$ws.Range("M17").Value = 15.87520469422147

$ws.Range("B18").Value = 8.081322329906122
$ws.Range("D18").Value = 9.5440612668129
$ws.Range("E18").Value = 11.66424438942017
$ws.Range("F18").Value = 53.08861179876224
$ws.Range("G18").Value = 3.682246781561322
$ws.Range("K18").Value = 12.82161278012936
$ws.Range("M18").Value = 15.83177015263796

$ws.Range("B19").Value = 8.074410586006847
$ws.Range("D19").Value = 9.53106560258562
$ws.Range("E19").Value = 11.64400458925795
$ws.Range("F19").Value = 52.97618221892879
$ws.Range("G19").Value = 3.682598216771344
$ws.Range("K19").Value = 12.80669546554483
$ws.Range("M19").Value = 15.81714112314018

$ws.Range("B20").Value = 8.105594205165403
$ws.Range("D20").Value = 9.589398895108896
$ws.Range("E20").Value = 11.73498124103717
$ws.Range("F20").Value = 53.48064229201559
$ws.Range("G20").Value = 3.681025487374313
$ws.Range("K20").Value = 12.87419101665644
$ws.Range("M20").Value = 15.88327977239874

$ws.Range("B21").Value = 8.211679956835635
$ws.Range("D21").Value = 9.782701343622477
$ws.Range("E21").Value = 12.03855513752806
$ws.Range("F21").Value = 55.14833127842075
$ws.Range("G21").Value = 3.675894134293604
$ws.Range("K21").Value = 13.10710439626872
$ws.Range("M21").Value = 16.11066211115691

$ws.Range("B22").Value = 8.281886066224351
$ws.Range("D22").Value = 9.906927473946199
$ws.Range("E22").Value = 12.23508853701038
$ws.Range("F22").Value = 56.21655476998905
$ws.Range("G22").Value = 3.672653606608187
$ws.Range("K22").Value = 13.26357414829513
$ws.Range("M22").Value = 16.26287831143585

$ws.Range("B23").Value = 8.244342667731715
$ws.Range("D23").Value = 9.840825040810033
$ws.Range("E23").Value = 12.13038524687753
$ws.Range("F23").Value = 55.64849675192235
$ws.Range("G23").Value = 3.674372837770516
$ws.Range("K23").Value = 13.17969580296093
$ws.Range("M23").Value = 16.18132311641847

$ws.Range("B24").Value = 8.103876541421148
$ws.Range("D24").Value = 9.586205355933648
$ws.Range("E24").Value = 11.7299923119449
$ws.Range("F24").Value = 53.45303841283113
$ws.Range("G24").Value = 3.681111277333899
$ws.Range("K24").Value = 12.87046056831439
$ws.Range("M24").Value = 15.87962771390743

$ws.Range("B25").Value = 7.956791470951344
$ws.Range("D25").Value = 9.303189928266157
$ws.Range("E25").Value = 11.2921154313787
$ws.Range("F25").Value = 51.00117157631593
$ws.Range("G25").Value = 3.688870448299749
$ws.Range("K25").Value = 12.55723545469825
$ws.Range("M25").Value = 15.57123934598693
